# Regenerate save_data "K" column (column G) values.
# The underlying change replaces the old "Strike#" derived K values with
# freshly calculated K values (std/mean based s_vals calc), row by row,
# for the single data sheet in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K (column G) value, as produced by the
# regeneration of save_data (Strike# -> K) and the std/mean / s_vals calc.
$kUpdates = [ordered]@{
    2  = 2
    3  = 1
    4  = 2
    5  = 0
    6  = 1
    7  = 0
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 0
    17 = 4
    18 = 1
    20 = 2
    21 = 1
    22 = 1
    23 = 3
    24 = 1
}

foreach ($row in $kUpdates.Keys) {
    $ws.Range("G$row").Value = $kUpdates[$row]
}

$wb.Save()
